$d = $word.ActiveDocument

# --- "Programa resumido" (PT) paragraph: split into 4 lines with manual line breaks ---
$d.Content.Find.Execute("Qualidade2 – Melhoramentos", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Qualidade^l2 – Melhoramentos", 2)
$d.Content.Find.Execute("Produção3 – Desafios", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Produção^l3 – Desafios", 2)
$d.Content.Find.Execute("produção4 – Controle", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "produção^l4 – Controle", 2)

# --- "Programa resumido" (EN, italic) paragraph: split into 4 lines ---
$d.Content.Find.Execute("Control2 - Production Improvements", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Control^l2 - Production Improvements", 2)
$d.Content.Find.Execute("Improvements3 - Production challenges", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Improvements^l3 - Production challenges", 2)
$d.Content.Find.Execute("challenges4 - Quality Control", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "challenges^l4 - Quality Control", 2)

# --- "Programa" (PT) paragraph: split into 8 lines ---
$d.Content.Find.Execute("QualidadeIntrodução. Planejamento", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Qualidade^lIntrodução. Planejamento", 2)
$d.Content.Find.Execute("qualidade.2 – Melhoramentos", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "qualidade.^l2 – Melhoramentos", 2)
$d.Content.Find.Execute("ProduçãoIntrodução. Medidas", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Produção^lIntrodução. Medidas", 2)
$d.Content.Find.Execute("Total.3 – Desafios", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Total.^l3 – Desafios", 2)
$d.Content.Find.Execute("produçãoIntrodução. Tipo", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "produção^lIntrodução. Tipo", 2)
$d.Content.Find.Execute("estratégias.4 - CONTROLE DA QUALIDADE", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "estratégias.^l4 - CONTROLE DA QUALIDADE", 2)
$d.Content.Find.Execute("QUALIDADEAs Sete Ferramentas", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "QUALIDADE^lAs Sete Ferramentas", 2)

# --- "Programa" (EN, italic) paragraph: split into 8 lines ---
$d.Content.Find.Execute("ControlIntroduction. Planning Quality Control.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Control^lIntroduction. Planning Quality Control.", 2)
$d.Content.Find.Execute("Control.2 - Production Improvements", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Control.^l2 - Production Improvements", 2)
$d.Content.Find.Execute("ImprovementsIntroduction. Measures", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Improvements^lIntroduction. Measures", 2)
$d.Content.Find.Execute("Management.3 - Production challenges", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Management.^l3 - Production challenges", 2)
$d.Content.Find.Execute("challengesIntroduction. Types", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "challenges^lIntroduction. Types", 2)
$d.Content.Find.Execute("strategies.4 - QUALITY CONTROL", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "strategies.^l4 - QUALITY CONTROL", 2)
$d.Content.Find.Execute("CONTROLThe Seven Quality Tools", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CONTROL^lThe Seven Quality Tools", 2)

# --- "Bibliografia" paragraph: split into 3 lines with double manual line breaks ---
$d.Content.Find.Execute("2002. VENANZI", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2002. ^l^lVENANZI", 2)
$d.Content.Find.Execute("2014Textos complementares", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2014^l^lTextos complementares", 2)
